$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Slide 3: remove the "A verifier" textbox ("ZoneTexte 2", id=3)
#    bottom-right reminder callout.
# ------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
for ($i = $s3.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s3.Shapes.Item($i)
    if ($shape.Name -eq "ZoneTexte 2" -and $shape.Id -eq 3) {
        $shape.Delete()
    }
}

# ------------------------------------------------------------------
# 2) Slide 4: fix the "Occuper : (smiley) / Chauffer : (fire)" text
#    -> "Occupe : (smiley) / Chauffe : (fire)" merged into a single
#    run, and shrink the autofit textbox to its new width.
# ------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
for ($i = 1; $i -le $s4.Shapes.Count; $i++) {
    $shape = $s4.Shapes.Item($i)
    if ($shape.Name -eq "ZoneTexte 8" -and $shape.Id -eq 9) {
        $tr = $shape.TextFrame.TextRange

        # Drop the first run ("Occuper : (smiley) ") so that re-typing
        # the whole text collapses onto the run that is left (which
        # already carries dirty="0"), instead of keeping two runs.
        $firstRun = $tr.Characters(1, 12)
        $firstRun.Delete()

        $tr2 = $shape.TextFrame.TextRange
        $newText = "Occup" + [char]0x00E9 + " : " + [char]0x1F642 + " / Chauff" + [char]0x00E9 + " : " + [char]0x1F525 + "  "
        $tr2.Text = $newText

        # spAutoFit text box: width shrinks to match the corrected text.
        $shape.Width = 151.71733
    }
}
